$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "https://www.booklooker.de/B%C3%BCcher/Angebote/sparteID=581&searchUserTyp=2&zustand=2&paymentMethods%5B6%5D=1?sortOrder=preis_total"
